$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 75; this pushes the existing rows 75-169
# down to 76-170 (along with their formatting), matching the diff which
# shows every row from 76 to 170 taking on the values that used to sit
# one row above it.
$ws.Rows("75:75").Insert()

# Fill in the brand-new row 75 with its own data (not shifted from
# anywhere else).
$ws.Range("A75").Value = 4
$ws.Range("B75").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C75").Value = "Los Lagos"
$ws.Range("D75").Value = 45174
$ws.Range("D75").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E75").Value = 10
$ws.Range("F75").Value = 100112031
$ws.Range("G75").Value = "Poroto verde"
$ws.Range("H75").Value = "Magnum"
$ws.Range("I75").Value = "Primera"
$ws.Range("J75").Value = 45
$ws.Range("K75").Value = 30000
$ws.Range("L75").Value = 30000
$ws.Range("M75").Value = 30000
$ws.Range("N75").Value = "$/saco 25 kilos"
$ws.Range("O75").Value = "Perú"
$ws.Range("P75").Value = 1200
$ws.Range("Q75").Value = 25
$ws.Range("R75").Value = "Hortaliza"
